# Rebuilt school networks table: apply the "Compact" paragraph style to
# the empty <w:p/> paragraphs in the "Due" column cells (3rd column) of
# the first table, matching the non-empty cells in the same table/column.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$rowCount = $t.Rows.Count
$colCount = $t.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $t.Cell($r, $colCount)
    $para = $cell.Range.Paragraphs.Item(1)

    # Cell.Range.Text includes the paragraph mark (CR) and cell-end mark
    # (BEL) even when the cell has no visible content, so strip those
    # control characters before testing for emptiness.
    $visibleText = $cell.Range.Text -replace "[\x07\x0d]", ""

    if ($visibleText -eq "" -and $para.Style.NameLocal -ne "Compact") {
        $para.Style = "Compact"
    }
}
